$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.513.03"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.418.21"
$ws.Range("E3").Value = "  +8.92%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.41"
$ws.Range("E5").Value = "  +11.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.98"
$ws.Range("E6").Value = "  -4.01%  "
$ws.Range("E7").Value = "  +3.99%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.661"
$ws.Range("E9").Value = "  +11.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.29"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.05"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.52"
$ws.Range("E14").Value = "  +18.48%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.781.65"
$ws.Range("E16").Value = "  +8.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.410.33"
$ws.Range("E17").Value = "  +7.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.559.75"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +6.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  +6.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.67"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "263.05"
$ws.Range("E23").Value = "  +14.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.76"
$ws.Range("E25").Value = "  +10.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("E26").Value = "  +6.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.97"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.85"
$ws.Range("E29").Value = "  +10.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "178.82"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.06"
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0935"
$ws.Range("E34").Value = "  +7.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("E35").Value = "  +7.54%  "
$ws.Range("E36").Value = "  +7.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.90"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0373"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("E41").Value = "  +22.93%  "
$ws.Range("E42").Value = "  +26.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.03"
$ws.Range("E43").Value = "  +26.61%  "
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.36"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.67"
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.72"
$ws.Range("E48").Value = "  +16.39%  "
$ws.Range("E49").Value = "  +7.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.33"
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.611.78"
$ws.Range("E51").Value = "  +14.80%  "
